# [Word] Add ChangeTracking snippet
# Adds three new rows to the "Snippets" table on the "Snippets" sheet:
#   - Document / changeTrackingMode / (blank) / word-manage-change-tracking / getChangeTrackingMode
#   - Document / changeTrackingMode / (blank) / word-manage-change-tracking / setChangeTrackingMode
#   - Range    / getReviewedText    / 1        / word-manage-change-tracking / getReviewedText

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Insert the two new "Document / changeTrackingMode" rows right above the
# --- existing "Document / properties" row (old row 22), pushing everything
# --- below it down by two rows.
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()

# --- Insert the new "Range / getReviewedText" row right above the existing
# --- "Range / getTextRanges" row (old row 45, now at row 47 after the two
# --- inserts above), pushing everything below it down by one more row.
$ws.Rows.Item(47).Insert()

# --- Copy formatting onto the two new Document rows from a row that already
# --- carries the same visual style used elsewhere in the "Range"/"Document"
# --- blocks (row 46 here == old row 44, "Range / getComments").
$ws.Range("A46:E46").Copy()
$ws.Range("A22:E22").PasteSpecial(-4122)
$ws.Range("A46:E46").Copy()
$ws.Range("A23:E23").PasteSpecial(-4122)

# --- Row 47 (the new Range/getReviewedText row) already inherited the
# --- correct formatting automatically from the row above it on insert.

# --- Fill in the new row values.
$ws.Range("A22").Value = "Document"
$ws.Range("B22").Value = "changeTrackingMode"
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = "word-manage-change-tracking"
$ws.Range("E22").Value = "getChangeTrackingMode"

$ws.Range("A23").Value = "Document"
$ws.Range("B23").Value = "changeTrackingMode"
$ws.Range("C23").Value = ""
$ws.Range("D23").Value = "word-manage-change-tracking"
$ws.Range("E23").Value = "setChangeTrackingMode"

$ws.Range("A47").Value = "Range"
$ws.Range("B47").Value = "getReviewedText"
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = "word-manage-change-tracking"
$ws.Range("E47").Value = "getReviewedText"

# --- Grow the table (ListObject) so it covers the three new rows, and keep
# --- the autofilter/sort-state ranges in sync with it.
$lo.Resize($ws.Range("A1:E55"))

# --- Match the saved sheet view (frozen pane / selection) from the edit.
$sheetView = $ws.Application
$ws.Range("A22").Select()
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("A48").Select()
